$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs whose B:AC (columns 2..29) content needs to be swapped between
# the two rows (column A, the plain row index, stays untouched).
$pairs = @(
    @(297, 298),
    @(302, 303),
    @(309, 310),
    @(323, 325)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range($ws.Cells.Item($r1, 2), $ws.Cells.Item($r1, 29))
    $range2 = $ws.Range($ws.Cells.Item($r2, 2), $ws.Cells.Item($r2, 29))

    $temp = $range1.Value2
    $range1.Value2 = $range2.Value2
    $range2.Value2 = $temp
}
